$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 27; $r -le 170; $r++) {
    $ws.Range("H$r").Value = 0
}
